# Fruta / hortaliza, semanal
# Update Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O),
# Precio promedio ponderado (P) and Precio $/Kg (S) for rows 2-13.
# The new values are a row-permutation of the original ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @(45076, 20, 15000, 15000, 15000, 833)
    3  = @(45112, 20, 20000, 20000, 20000, 1111)
    4  = @(45072, 30, 15000, 15000, 15000, 833)
    5  = @(45092, 60, 18000, 19000, 18667, 1037)
    6  = @(45083, 50, 15000, 15000, 15000, 833)
    7  = @(45085, 30, 19000, 19000, 19000, 1056)
    8  = @(45096, 30, 20000, 20000, 20000, 1111)
    9  = @(45055, 50, 15000, 15000, 15000, 833)
    10 = @(45111, 20, 20000, 20000, 20000, 1111)
    11 = @(45084, 50, 18000, 19000, 18500, 1028)
    12 = @(45061, 40, 15000, 15000, 15000, 833)
    13 = @(45069, 60, 15000, 15000, 15000, 833)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value  = $vals[0]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $vals[1]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $vals[2]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $vals[3]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $vals[4]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $vals[5]   # S - Precio $/Kg
}
